$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Última fecha"
$ws.Range("D2").Value = "ANDINA"
$ws.Range("G2").Value = "2023-12-29 00:00:00"
$ws.Range("D3").Value = "ANDINA"
$ws.Range("F3").Value = "Sin precio"
$ws.Range("G3").Value = "2023-12-29 00:00:00"
